$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric identifiers
$ws.Range("A2").Value = 80139019
$ws.Range("B2").Value = 88921
$ws.Range("E2").Value = 5741

# Species name / scientific name / author
$ws.Range("F2").Value = "Tjockfotad fingersvamp"
$ws.Range("G2").Value = "Ramaria flavescens"
$ws.Range("H2").Value = "(Schaeff.) R. H. Petersen"

# Antal (I2): cleared to an empty (but still present/text-typed) cell.
# A bare quote-prefix forces an empty Text cell instead of deleting it,
# then resetting the style avoids leaving a stray quote-prefix format.
$ci2 = $ws.Range("I2")
$ci2.Value = "'"
$ci2.Style = "Normal"

# Enhet (J2), Ålder-Stadium (K2), Kön (L2), Metod (N2): fully removed
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("N2").ClearContents()

# Locality name
$ws.Range("P2").Value = "Tveta friluftsgård, 300 m V om, Srm"

# Coordinates
$ws.Range("Q2").Value = 648222.682956806
$ws.Range("R2").Value = 6560420.292955686
$ws.Range("S2").Value = 50

# Dates (stored as plain text, not real Excel dates): force text format
# before assigning so Excel doesn't auto-convert the string to a serial
# date, then reset the style back so no stray number-format sticks around.
$cy2 = $ws.Range("Y2")
$cy2.NumberFormat = "@"
$cy2.Value = "2019-09-27"
$cy2.Style = "Normal"

$caa2 = $ws.Range("AA2")
$caa2.NumberFormat = "@"
$caa2.Value = "2019-09-27"
$caa2.Style = "Normal"

# Bestämningsmetod (AF2): fully removed
$ws.Range("AF2").ClearContents()

# Biotop-beskrivning (AI2): newly added
$ws.Range("AI2").Value = "barrskog"

# Observer names
$ws.Range("AW2").Value = "Hans Rydberg"
$ws.Range("AX2").Value = "Hans Rydberg"
